# "Generate Report for Archive"
#
# 1) Every "Ready for handoff" status cell becomes "In Translation"
#    (Overview!E2:F2, E3:F3 ; zh-cn!C2:C3 ; de-de!C2:C3).
# 2) The (now narrower) status columns shrink:
#    Overview columns E & F, and column C on both the zh-cn and de-de
#    sheets, go from ~17.22 chars wide to ~13.41 chars wide.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2").Value = "In Translation"
$ovw.Range("F2").Value = "In Translation"
$ovw.Range("E3").Value = "In Translation"
$ovw.Range("F3").Value = "In Translation"

# Target stored width is 13.4101845877511 character-units; the COM
# ColumnWidth setter snaps to the sheet's pixel grid, and 12.5 is the
# closest input that lands on that value.
$ovw.Columns.Item(5).ColumnWidth = 12.5
$ovw.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet ------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.5
